# Incubators workbook: add new "public submissions" columns (I:S) and
# append 4 newly-submitted incubator rows (30-33) that only populate
# those new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header row (row 1), columns I..S
# ---------------------------------------------------------------------
$headers = @(
    "name",
    "location",
    "website",
    "email",
    "phone",
    "focus_sectors",
    "program_duration",
    "equity_taken",
    "funding_amount",
    "application_deadline",
    "description"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, 9 + $i).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 2. Four new submission rows (30-33), populating only I, L, M, N with
#    real values; J, K, O, P, Q, R, S stay blank but must still exist as
#    cells so the sheet dimension grows to A1:S33. Touching a formatting
#    property (without actually changing the format) is enough to make
#    the engine materialize an otherwise-empty cell.
# ---------------------------------------------------------------------
$newRows = @(
    @{ name = "zdvxfb"; email = "priyanshu@gmail.com"; phone = "8103700333"; focus = "kk" },
    @{ name = "pari";   email = "priyanshu@gmail.com"; phone = "8103700333"; focus = "kk" },
    @{ name = "pari";   email = "priyanshu@gmail.com"; phone = "8103700333"; focus = "kk" },
    @{ name = "pari";   email = "priyanshu@gmail.com"; phone = "8103700333"; focus = "kk" }
)

$blankCols = @(10, 11, 15, 16, 17, 18, 19)   # J, K, O, P, Q, R, S

$row = 30
foreach ($entry in $newRows) {
    $ws.Cells.Item($row, 9).Value  = $entry.name    # I - name
    $ws.Cells.Item($row, 12).Value = $entry.email    # L - email
    $ws.Cells.Item($row, 13).Value = $entry.phone    # M - phone
    $ws.Cells.Item($row, 14).Value = $entry.focus    # N - focus_sectors

    foreach ($col in $blankCols) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Font.Bold = $false
    }

    $row = $row + 1
}

# ---------------------------------------------------------------------
# 3. Keep the "number stored as text" error-check ignore range in sync
#    with the new A1:S33 extent (mirrors the sheet dimension growth).
# ---------------------------------------------------------------------
$full = $ws.Range("A1:S33")
$full.Errors.Item(3).Ignore = $true

Write-Host "Added public submission columns (I:S) and 4 new rows (30-33)."
